# Auto-generated from the authoritative diff: updates cached market-board
# price/profit figures (columns H-N) across 42 leve rows spanning all 8
# crafter sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# ALC row 9
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 639580.9399999999
$ws.Range("I9").Value = 1050602.9
$ws.Range("K9").Value = 1050602.9
$ws.Range("M9").Value = -1050433.9

# ALC row 17
$ws.Range("H17").Value = 1913.8334
$ws.Range("J17").Value = 2028.125
$ws.Range("L17").Value = 6084.375
$ws.Range("N17").Value = -6420.375

# ALC row 62
$ws.Range("H62").Value = 4685.1816
$ws.Range("I62").Value = 4637.2
$ws.Range("J62").Value = 4725.1665
$ws.Range("K62").Value = 4637.2
$ws.Range("L62").Value = 4725.1665
$ws.Range("M62").Value = -4013.2
$ws.Range("N62").Value = -5973.1665

# ALC row 65
$ws.Range("H65").Value = 4685.1816
$ws.Range("I65").Value = 4637.2
$ws.Range("J65").Value = 4725.1665
$ws.Range("K65").Value = 23186
$ws.Range("L65").Value = 23625.8325
$ws.Range("M65").Value = -20066
$ws.Range("N65").Value = -29865.8325

# ALC row 100
$ws.Range("H100").Value = 1166.8334
$ws.Range("I100").Value = 1166.8334
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 1166.8334
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -625.8334
$ws.Range("N100").ClearContents()

# ALC row 132
$ws.Range("H132").Value = 5023.9688
$ws.Range("I132").Value = 5163.5806
$ws.Range("J132").Value = 696
$ws.Range("K132").Value = 15490.7418
$ws.Range("L132").Value = 2088
$ws.Range("M132").Value = -12960.7418
$ws.Range("N132").Value = -7148

# ALC row 138
$ws.Range("H138").Value = 4069.2666
$ws.Range("I138").Value = 4849.2
$ws.Range("J138").Value = 3846.4285
$ws.Range("K138").Value = 14547.6
$ws.Range("L138").Value = 11539.2855
$ws.Range("M138").Value = -9407.599999999999
$ws.Range("N138").Value = -21819.2855

# ARM row 97
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 1196.7667
$ws.Range("J97").Value = 1123.875
$ws.Range("L97").Value = 1123.875
$ws.Range("N97").Value = -2115.875

# ARM row 110
$ws.Range("H110").Value = 765.7143
$ws.Range("J110").Value = 812
$ws.Range("L110").Value = 812
$ws.Range("N110").Value = -4902

# ARM row 132
$ws.Range("H132").Value = 4220.3413
$ws.Range("I132").Value = 3045.2727
$ws.Range("J132").Value = 5580.9473
$ws.Range("K132").Value = 9135.8181
$ws.Range("L132").Value = 16742.8419
$ws.Range("M132").Value = -6605.8181
$ws.Range("N132").Value = -21802.8419

# BSM row 20
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 14289057
$ws.Range("J20").Value = 2534
$ws.Range("L20").Value = 2534
$ws.Range("N20").Value = -3028

# BSM row 105
$ws.Range("H105").Value = 10002099
$ws.Range("I105").Value = 771393.3
$ws.Range("J105").Value = 19232804
$ws.Range("K105").Value = 771393.3
$ws.Range("L105").Value = 19232804
$ws.Range("M105").Value = -769646.3
$ws.Range("N105").Value = -19236298

# BSM row 134
$ws.Range("H134").Value = 3142.9
$ws.Range("I134").Value = 2796.889
$ws.Range("J134").Value = 6257
$ws.Range("K134").Value = 8390.667000000001
$ws.Range("L134").Value = 18771
$ws.Range("M134").Value = -5855.667000000001
$ws.Range("N134").Value = -23841

# CRP row 59
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H59").Value = 69463.5
$ws.Range("J59").Value = 69463.5
$ws.Range("L59").Value = 69463.5
$ws.Range("N59").Value = -71753.5

# CRP row 74
$ws.Range("H74").Value = 68556.75
$ws.Range("J74").Value = 81409
$ws.Range("L74").Value = 81409
$ws.Range("N74").Value = -83157

# CRP row 77
$ws.Range("H77").Value = 68556.75
$ws.Range("J77").Value = 81409
$ws.Range("L77").Value = 244227
$ws.Range("N77").Value = -252963

# CRP row 86
$ws.Range("H86").Value = 7563.857
$ws.Range("I86").Value = 7491.1665
$ws.Range("K86").Value = 7491.1665
$ws.Range("M86").Value = -6368.1665

# CRP row 89
$ws.Range("H89").Value = 7563.857
$ws.Range("I89").Value = 7491.1665
$ws.Range("K89").Value = 37455.8325
$ws.Range("M89").Value = -31839.8325

# CRP row 122
$ws.Range("H122").Value = 3050.8696
$ws.Range("I122").Value = 3108.7
$ws.Range("K122").Value = 9326.099999999999
$ws.Range("M122").Value = -6876.099999999999

# CRP row 134
$ws.Range("H134").Value = 5450.727
$ws.Range("I134").Value = 4856.8184
$ws.Range("K134").Value = 14570.4552
$ws.Range("M134").Value = -12035.4552

# CUL row 59
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H59").Value = 4855.7144
$ws.Range("I59").Value = 990
$ws.Range("K59").Value = 2970
$ws.Range("M59").Value = -2430

# CUL row 122
$ws.Range("H122").Value = 2034.25
$ws.Range("J122").Value = 2054.7222
$ws.Range("L122").Value = 18492.4998
$ws.Range("N122").Value = -23392.4998

# CUL row 132
$ws.Range("H132").Value = 2094.8
$ws.Range("I132").Value = 2490
$ws.Range("J132").Value = 1996
$ws.Range("K132").Value = 22410
$ws.Range("L132").Value = 17964
$ws.Range("M132").Value = -19880
$ws.Range("N132").Value = -23024

# CUL row 134
$ws.Range("H134").Value = 2250.5
$ws.Range("I134").Value = 2250.5
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 6751.5
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -1681.5
$ws.Range("N134").ClearContents()

# GSM row 70
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 73707.75999999999
$ws.Range("I70").Value = 91423.13
$ws.Range("K70").Value = 91423.13
$ws.Range("M70").Value = -91153.13

# GSM row 73
$ws.Range("H73").Value = 73707.75999999999
$ws.Range("I73").Value = 91423.13
$ws.Range("K73").Value = 91423.13
$ws.Range("M73").Value = -90487.13

# GSM row 80
$ws.Range("H80").Value = 38464056
$ws.Range("I80").Value = 100001450
$ws.Range("J80").Value = 3185.3125
$ws.Range("K80").Value = 100001450
$ws.Range("L80").Value = 3185.3125
$ws.Range("M80").Value = -100000452
$ws.Range("N80").Value = -5181.3125

# GSM row 83
$ws.Range("H83").Value = 38464056
$ws.Range("I83").Value = 100001450
$ws.Range("J83").Value = 3185.3125
$ws.Range("K83").Value = 500007250
$ws.Range("L83").Value = 15926.5625
$ws.Range("M83").Value = -500002258
$ws.Range("N83").Value = -25910.5625

# GSM row 97
$ws.Range("H97").Value = 1935.6666
$ws.Range("I97").Value = 1826.4667
$ws.Range("K97").Value = 1826.4667
$ws.Range("M97").Value = -1330.4667

# GSM row 98
$ws.Range("H98").Value = 24500
$ws.Range("J98").Value = 24500
$ws.Range("L98").Value = 24500
$ws.Range("N98").Value = -30490

# GSM row 122
$ws.Range("H122").Value = 3250.9443
$ws.Range("I122").Value = 2988.5
$ws.Range("J122").Value = 4563.1665
$ws.Range("K122").Value = 8965.5
$ws.Range("L122").Value = 13689.4995
$ws.Range("M122").Value = -6515.5
$ws.Range("N122").Value = -18589.4995

# GSM row 126
$ws.Range("H126").Value = 14832.111
$ws.Range("I126").Value = 6498
$ws.Range("J126").Value = 23166.223
$ws.Range("K126").Value = 19494
$ws.Range("L126").Value = 69498.66900000001
$ws.Range("M126").Value = -17024
$ws.Range("N126").Value = -74438.66900000001

# LTW row 7
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5290.6
$ws.Range("J7").Value = 7781
$ws.Range("L7").Value = 7781
$ws.Range("N7").Value = -8005

# LTW row 16
$ws.Range("H16").Value = 711.069
$ws.Range("I16").Value = 736.3333
$ws.Range("J16").Value = 589.8
$ws.Range("K16").Value = 736.3333
$ws.Range("L16").Value = 589.8
$ws.Range("M16").Value = -566.3333
$ws.Range("N16").Value = -929.8

# LTW row 61
$ws.Range("H61").Value = 1965.4546
$ws.Range("I61").Value = 2119.5
$ws.Range("J61").Value = 425
$ws.Range("K61").Value = 2119.5
$ws.Range("L61").Value = 425
$ws.Range("M61").Value = -1917.5
$ws.Range("N61").Value = -829

# LTW row 100
$ws.Range("H100").Value = 4173.136
$ws.Range("I100").Value = 3042.6316
$ws.Range("J100").Value = 11333
$ws.Range("K100").Value = 3042.6316
$ws.Range("L100").Value = 11333
$ws.Range("M100").Value = -2501.6316
$ws.Range("N100").Value = -12415

# LTW row 113
$ws.Range("H113").Value = 1965.4546
$ws.Range("I113").Value = 2119.5
$ws.Range("J113").Value = 425
$ws.Range("K113").Value = 2119.5
$ws.Range("L113").Value = 425
$ws.Range("M113").Value = 50.5
$ws.Range("N113").Value = -4765

# LTW row 122
$ws.Range("H122").Value = 6420.885
$ws.Range("I122").Value = 5312.05
$ws.Range("J122").Value = 10117
$ws.Range("K122").Value = 15936.15
$ws.Range("L122").Value = 30351
$ws.Range("M122").Value = -13486.15
$ws.Range("N122").Value = -35251

# LTW row 126
$ws.Range("H126").Value = 5290.6
$ws.Range("J126").Value = 7781
$ws.Range("L126").Value = 23343
$ws.Range("N126").Value = -28283

# WVR row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 14707679
$ws.Range("I122").Value = 1983.5385
$ws.Range("J122").Value = 62501190
$ws.Range("K122").Value = 5950.6155
$ws.Range("L122").Value = 187503570
$ws.Range("M122").Value = -3500.6155
$ws.Range("N122").Value = -187508470

# WVR row 126
$ws.Range("H126").Value = 11273.091
$ws.Range("I126").Value = 13222.667
$ws.Range("J126").Value = 2500
$ws.Range("K126").Value = 39668.001
$ws.Range("L126").Value = 7500
$ws.Range("M126").Value = -37198.001
$ws.Range("N126").Value = -12440

# WVR row 132
$ws.Range("H132").Value = 2404.7
$ws.Range("I132").Value = 2007.4231
$ws.Range("J132").Value = 4987
$ws.Range("K132").Value = 6022.2693
$ws.Range("L132").Value = 14961
$ws.Range("M132").Value = -3492.2693
$ws.Range("N132").Value = -20021

